$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 25.06781033333333
$ws.Range("H2").Value = 75.20343099999999
$ws.Range("I2").Value = 0.7308832858982242
$ws.Range("J2").Value = 0.7308832858982242
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 34.58696433333333
$ws.Range("N2").Value = 103.760893
$ws.Range("O2").Value = 0.9801058085769325
$ws.Range("P2").Value = 0.9801058085769326
$ws.Range("Q2").Value = 867.0194619137648
$ws.Range("R2").Value = 7803.175157223882
$ws.Range("S2").Value = 0.7163429539006443
$ws.Range("T2").Value = 0.7163429539006444

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 25.06781033333333
$ws.Range("H3").Value = 75.20343099999999
$ws.Range("I3").Value = 0.7308832858982242
$ws.Range("J3").Value = 0.7308832858982242
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.2577526666666667
$ws.Range("N3").Value = 0.773258
$ws.Range("O3").Value = 0.007304049101896046
$ws.Range("P3").Value = 0.007304049101896047
$ws.Range("Q3").Value = 6.461294960910888
$ws.Range("R3").Value = 58.151654648198
$ws.Range("S3").Value = 0.005338407407955755
$ws.Range("T3").Value = 0.005338407407955756

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 25.06781033333333
$ws.Range("H4").Value = 75.20343099999999
$ws.Range("I4").Value = 0.7308832858982242
$ws.Range("J4").Value = 0.7308832858982242
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.4442936666666666
$ws.Range("N4").Value = 1.332881
$ws.Range("O4").Value = 0.01259014232117133
$ws.Range("P4").Value = 0.01259014232117133
$ws.Range("Q4").Value = 11.13746936830122
$ws.Range("R4").Value = 100.237224314711
$ws.Range("S4").Value = 0.009201924589623999
$ws.Range("T4").Value = 0.009201924589624001

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.8930513333333332
$ws.Range("H5").Value = 2.679154
$ws.Range("I5").Value = 0.02603802583086097
$ws.Range("J5").Value = 0.02603802583086097
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 34.58696433333333
$ws.Range("N5").Value = 103.760893
$ws.Range("O5").Value = 0.9801058085769325
$ws.Range("P5").Value = 0.9801058085769326
$ws.Range("Q5").Value = 30.88793461383577
$ws.Range("R5").Value = 277.9914115245219
$ws.Range("S5").Value = 0.02552002036070305
$ws.Range("T5").Value = 0.02552002036070305

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.8930513333333332
$ws.Range("H6").Value = 2.679154
$ws.Range("I6").Value = 0.02603802583086097
$ws.Range("J6").Value = 0.02603802583086097
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.2577526666666667
$ws.Range("N6").Value = 0.773258
$ws.Range("O6").Value = 0.007304049101896046
$ws.Range("P6").Value = 0.007304049101896047
$ws.Range("Q6").Value = 0.2301863626368889
$ws.Range("R6").Value = 2.071677263732
$ws.Range("S6").Value = 0.0001901830191850461
$ws.Range("T6").Value = 0.0001901830191850461

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.8930513333333332
$ws.Range("H7").Value = 2.679154
$ws.Range("I7").Value = 0.02603802583086097
$ws.Range("J7").Value = 0.02603802583086097
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.4442936666666666
$ws.Range("N7").Value = 1.332881
$ws.Range("O7").Value = 0.01259014232117133
$ws.Range("P7").Value = 0.01259014232117133
$ws.Range("Q7").Value = 0.3967770514082221
$ws.Range("R7").Value = 3.570993462673999
$ws.Range("S7").Value = 0.0003278224509728751
$ws.Range("T7").Value = 0.0003278224509728751

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.337104666666667
$ws.Range("H8").Value = 25.011314
$ws.Range("I8").Value = 0.2430786882709149
$ws.Range("J8").Value = 0.2430786882709149
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 34.58696433333333
$ws.Range("N8").Value = 103.760893
$ws.Range("O8").Value = 0.9801058085769325
$ws.Range("P8").Value = 0.9801058085769326
$ws.Range("Q8").Value = 288.3551417492669
$ws.Range("R8").Value = 2595.196275743402
$ws.Range("S8").Value = 0.2382428343155852
$ws.Range("T8").Value = 0.2382428343155852

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.337104666666667
$ws.Range("H9").Value = 25.011314
$ws.Range("I9").Value = 0.2430786882709149
$ws.Range("J9").Value = 0.2430786882709149
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.2577526666666667
$ws.Range("N9").Value = 0.773258
$ws.Range("O9").Value = 0.007304049101896046
$ws.Range("P9").Value = 0.007304049101896047
$ws.Range("Q9").Value = 2.148910960112445
$ws.Range("R9").Value = 19.340198641012
$ws.Range("S9").Value = 0.001775458674755245
$ws.Range("T9").Value = 0.001775458674755245

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.337104666666667
$ws.Range("H10").Value = 25.011314
$ws.Range("I10").Value = 0.2430786882709149
$ws.Range("J10").Value = 0.2430786882709149
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.4442936666666666
$ws.Range("N10").Value = 1.332881
$ws.Range("O10").Value = 0.01259014232117133
$ws.Range("P10").Value = 0.01259014232117133
$ws.Range("Q10").Value = 3.704122801737111
$ws.Range("R10").Value = 33.337105215634
$ws.Range("S10").Value = 0.00306039528057446
$ws.Range("T10").Value = 0.00306039528057446

